$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- 1) Bump the global 'Forandrad' (last-changed) date column C for every data row (2-113) ----
for ($r = 2; $r -le 113; $r++) {
    $ws.Cells.Item($r, 3).Value = 46063
}

# ---- 2) Rewrite rows whose content moved to a new row position ----
# Each block below writes columns A,B,G:Q,R for one destination row, sourced from the
# pre-edit content of the case that now belongs there (full-file row re-sort).

# Row 3  <=  case 'A 39165-2024'
$ws.Cells.Item(3, 1).Value = 'A 39165-2024'
$ws.Cells.Item(3, 2).Value = 45548
$rngGQ = $ws.Range("G3:Q3")
$rngGQ.Value = @(2.9,0,0,0,0,1,1,0,2,2,2)
$ws.Cells.Item(3, 18).Value = 'Skogsalm`r`nAsk'

# Row 4  <=  case 'A 38522-2025'
$ws.Cells.Item(4, 1).Value = 'A 38522-2025'
$ws.Cells.Item(4, 2).Value = 45884
$rngGQ = $ws.Range("G4:Q4")
$rngGQ.Value = @(2.6,1,2,0,0,0,0,0,0,0,2)
$ws.Cells.Item(4, 18).Value = 'Skogsknipprot`r`nStor häxört'

# Row 5  <=  case 'A 61064-2024'
$ws.Cells.Item(5, 1).Value = 'A 61064-2024'
$ws.Cells.Item(5, 2).Value = 45645.49443287037
$rngGQ = $ws.Range("G5:Q5")
$rngGQ.Value = @(5.6,1,1,0,0,1,0,0,1,1,2)
$ws.Cells.Item(5, 18).Value = 'Ask`r`nSkogsknipprot'

# Row 6  <=  case 'A 26855-2022'
$ws.Cells.Item(6, 1).Value = 'A 26855-2022'
$ws.Cells.Item(6, 2).Value = 44740
$rngGQ = $ws.Range("G6:Q6")
$rngGQ.Value = @(5.2,1,1,0,0,0,0,0,0,0,2)
$ws.Cells.Item(6, 18).Value = 'Skogsbräsma`r`nMattlummer'

# Row 7  <=  case 'A 15456-2024'
$ws.Cells.Item(7, 1).Value = 'A 15456-2024'
$ws.Cells.Item(7, 2).Value = 45401
$rngGQ = $ws.Range("G7:Q7")
$rngGQ.Value = @(2.5,0,0,0,1,0,0,0,1,1,1)
$ws.Cells.Item(7, 18).Value = 'Lundticka'

# Row 8  <=  case 'A 15475-2024'
$ws.Cells.Item(8, 1).Value = 'A 15475-2024'
$ws.Cells.Item(8, 2).Value = 45401
$rngGQ = $ws.Range("G8:Q8")
$rngGQ.Value = @(4.7,0,1,0,0,0,0,0,0,0,1)
$ws.Cells.Item(8, 18).Value = 'Strutbräken'

# Row 9  <=  case 'A 54424-2023'
$ws.Cells.Item(9, 1).Value = 'A 54424-2023'
$ws.Cells.Item(9, 2).Value = 45233
$rngGQ = $ws.Range("G9:Q9")
$rngGQ.Value = @(5.7,0,1,0,0,0,0,0,0,0,1)
$ws.Cells.Item(9, 18).Value = 'Myskmadra'

# Row 10  <=  case 'A 1577-2024'
$ws.Cells.Item(10, 1).Value = 'A 1577-2024'
$ws.Cells.Item(10, 2).Value = 45306
$rngGQ = $ws.Range("G10:Q10")
$rngGQ.Value = @(21.2,1,0,0,0,0,0,0,0,0,1)
$ws.Cells.Item(10, 18).Value = 'Större vattensalamander'

# Row 12  <=  case 'A 11170-2023'
$ws.Cells.Item(12, 1).Value = 'A 11170-2023'
$ws.Cells.Item(12, 2).Value = 44987
$rngGQ = $ws.Range("G12:Q12")
$rngGQ.Value = @(0.9,0,0,0,0,1,0,0,1,1,1)
$ws.Cells.Item(12, 18).Value = 'Korndådra'

# Row 25  <=  case 'A 13245-2024'
$ws.Cells.Item(25, 1).Value = 'A 13245-2024'
$ws.Cells.Item(25, 2).Value = 45386
$rngGQ = $ws.Range("G25:Q25")
$rngGQ.Value = @(13,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(25, 18).Value = ""

# Row 26  <=  case 'A 34859-2021'
$ws.Cells.Item(26, 1).Value = 'A 34859-2021'
$ws.Cells.Item(26, 2).Value = 44382
$rngGQ = $ws.Range("G26:Q26")
$rngGQ.Value = @(1.9,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(26, 18).Value = ""

# Row 27  <=  case 'A 21774-2024'
$ws.Cells.Item(27, 1).Value = 'A 21774-2024'
$ws.Cells.Item(27, 2).Value = 45442
$rngGQ = $ws.Range("G27:Q27")
$rngGQ.Value = @(0.6,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(27, 18).Value = ""

# Row 28  <=  case 'A 34810-2023'
$ws.Cells.Item(28, 1).Value = 'A 34810-2023'
$ws.Cells.Item(28, 2).Value = 45141
$rngGQ = $ws.Range("G28:Q28")
$rngGQ.Value = @(7,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(28, 18).Value = ""

# Row 29  <=  case 'A 22794-2024'
$ws.Cells.Item(29, 1).Value = 'A 22794-2024'
$ws.Cells.Item(29, 2).Value = 45448.49752314815
$rngGQ = $ws.Range("G29:Q29")
$rngGQ.Value = @(3.8,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(29, 18).Value = ""

# Row 30  <=  case 'A 35587-2024'
$ws.Cells.Item(30, 1).Value = 'A 35587-2024'
$ws.Cells.Item(30, 2).Value = 45531
$rngGQ = $ws.Range("G30:Q30")
$rngGQ.Value = @(1.6,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(30, 18).Value = ""

# Row 31  <=  case 'A 16667-2023'
$ws.Cells.Item(31, 1).Value = 'A 16667-2023'
$ws.Cells.Item(31, 2).Value = 45030
$rngGQ = $ws.Range("G31:Q31")
$rngGQ.Value = @(3.8,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(31, 18).Value = ""

# Row 32  <=  case 'A 38394-2025'
$ws.Cells.Item(32, 1).Value = 'A 38394-2025'
$ws.Cells.Item(32, 2).Value = 45883
$rngGQ = $ws.Range("G32:Q32")
$rngGQ.Value = @(0.5,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(32, 18).Value = ""

# Row 33  <=  case 'A 47696-2024'
$ws.Cells.Item(33, 1).Value = 'A 47696-2024'
$ws.Cells.Item(33, 2).Value = 45588.45855324074
$rngGQ = $ws.Range("G33:Q33")
$rngGQ.Value = @(3.3,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(33, 18).Value = ""

# Row 34  <=  case 'A 8332-2025'
$ws.Cells.Item(34, 1).Value = 'A 8332-2025'
$ws.Cells.Item(34, 2).Value = 45708.64819444445
$rngGQ = $ws.Range("G34:Q34")
$rngGQ.Value = @(0.5,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(34, 18).Value = ""

# Row 35  <=  case 'A 14543-2024'
$ws.Cells.Item(35, 1).Value = 'A 14543-2024'
$ws.Cells.Item(35, 2).Value = 45394
$rngGQ = $ws.Range("G35:Q35")
$rngGQ.Value = @(0.7,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(35, 18).Value = ""

# Row 36  <=  case 'A 22063-2025'
$ws.Cells.Item(36, 1).Value = 'A 22063-2025'
$ws.Cells.Item(36, 2).Value = 45785.37700231482
$rngGQ = $ws.Range("G36:Q36")
$rngGQ.Value = @(2,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(36, 18).Value = ""

# Row 37  <=  case 'A 20170-2024'
$ws.Cells.Item(37, 1).Value = 'A 20170-2024'
$ws.Cells.Item(37, 2).Value = 45434
$rngGQ = $ws.Range("G37:Q37")
$rngGQ.Value = @(0.6,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(37, 18).Value = ""

# Row 38  <=  case 'A 62403-2022'
$ws.Cells.Item(38, 1).Value = 'A 62403-2022'
$ws.Cells.Item(38, 2).Value = 44923
$rngGQ = $ws.Range("G38:Q38")
$rngGQ.Value = @(2.2,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(38, 18).Value = ""

# Row 39  <=  case 'A 19626-2021'
$ws.Cells.Item(39, 1).Value = 'A 19626-2021'
$ws.Cells.Item(39, 2).Value = 44309
$rngGQ = $ws.Range("G39:Q39")
$rngGQ.Value = @(17.8,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(39, 18).Value = ""

# Row 40  <=  case 'A 17980-2025'
$ws.Cells.Item(40, 1).Value = 'A 17980-2025'
$ws.Cells.Item(40, 2).Value = 45761.36854166666
$rngGQ = $ws.Range("G40:Q40")
$rngGQ.Value = @(3.4,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(40, 18).Value = ""

# Row 41  <=  case 'A 19282-2025'
$ws.Cells.Item(41, 1).Value = 'A 19282-2025'
$ws.Cells.Item(41, 2).Value = 45769
$rngGQ = $ws.Range("G41:Q41")
$rngGQ.Value = @(0.8,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(41, 18).Value = ""

# Row 42  <=  case 'A 15157-2023'
$ws.Cells.Item(42, 1).Value = 'A 15157-2023'
$ws.Cells.Item(42, 2).Value = 45016
$rngGQ = $ws.Range("G42:Q42")
$rngGQ.Value = @(2,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(42, 18).Value = ""

# Row 43  <=  case 'A 11491-2023'
$ws.Cells.Item(43, 1).Value = 'A 11491-2023'
$ws.Cells.Item(43, 2).Value = 44991
$rngGQ = $ws.Range("G43:Q43")
$rngGQ.Value = @(2.8,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(43, 18).Value = ""

# Row 44  <=  case 'A 42706-2025'
$ws.Cells.Item(44, 1).Value = 'A 42706-2025'
$ws.Cells.Item(44, 2).Value = 45908.37222222222
$rngGQ = $ws.Range("G44:Q44")
$rngGQ.Value = @(0.5,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(44, 18).Value = ""

# Row 45  <=  case 'A 42698-2025'
$ws.Cells.Item(45, 1).Value = 'A 42698-2025'
$ws.Cells.Item(45, 2).Value = 45908.3650462963
$rngGQ = $ws.Range("G45:Q45")
$rngGQ.Value = @(3.1,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(45, 18).Value = ""

# Row 46  <=  case 'A 43425-2025'
$ws.Cells.Item(46, 1).Value = 'A 43425-2025'
$ws.Cells.Item(46, 2).Value = 45911.4182175926
$rngGQ = $ws.Range("G46:Q46")
$rngGQ.Value = @(2.1,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(46, 18).Value = ""

# Row 47  <=  case 'A 29296-2021'
$ws.Cells.Item(47, 1).Value = 'A 29296-2021'
$ws.Cells.Item(47, 2).Value = 44361
$rngGQ = $ws.Range("G47:Q47")
$rngGQ.Value = @(0.7,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(47, 18).Value = ""

# Row 48  <=  case 'A 37907-2025'
$ws.Cells.Item(48, 1).Value = 'A 37907-2025'
$ws.Cells.Item(48, 2).Value = 45880
$rngGQ = $ws.Range("G48:Q48")
$rngGQ.Value = @(5.7,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(48, 18).Value = ""

# Row 49  <=  case 'A 10209-2023'
$ws.Cells.Item(49, 1).Value = 'A 10209-2023'
$ws.Cells.Item(49, 2).Value = 44986
$rngGQ = $ws.Range("G49:Q49")
$rngGQ.Value = @(1.2,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(49, 18).Value = ""

# Row 50  <=  case 'A 40470-2025'
$ws.Cells.Item(50, 1).Value = 'A 40470-2025'
$ws.Cells.Item(50, 2).Value = 45895
$rngGQ = $ws.Range("G50:Q50")
$rngGQ.Value = @(0.1,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(50, 18).Value = ""

# Row 51  <=  case 'A 44107-2023'
$ws.Cells.Item(51, 1).Value = 'A 44107-2023'
$ws.Cells.Item(51, 2).Value = 45188.43670138889
$rngGQ = $ws.Range("G51:Q51")
$rngGQ.Value = @(2.3,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(51, 18).Value = ""

# Row 52  <=  case 'A 47139-2025'
$ws.Cells.Item(52, 1).Value = 'A 47139-2025'
$ws.Cells.Item(52, 2).Value = 45929
$rngGQ = $ws.Range("G52:Q52")
$rngGQ.Value = @(1.1,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(52, 18).Value = ""

# Row 53  <=  case 'A 20235-2025'
$ws.Cells.Item(53, 1).Value = 'A 20235-2025'
$ws.Cells.Item(53, 2).Value = 45772.67386574074
$rngGQ = $ws.Range("G53:Q53")
$rngGQ.Value = @(0.9,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(53, 18).Value = ""

# Row 54  <=  case 'A 22262-2023'
$ws.Cells.Item(54, 1).Value = 'A 22262-2023'
$ws.Cells.Item(54, 2).Value = 45070
$rngGQ = $ws.Range("G54:Q54")
$rngGQ.Value = @(1.1,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(54, 18).Value = ""

# Row 55  <=  case 'A 48279-2021'
$ws.Cells.Item(55, 1).Value = 'A 48279-2021'
$ws.Cells.Item(55, 2).Value = 44449
$rngGQ = $ws.Range("G55:Q55")
$rngGQ.Value = @(1.3,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(55, 18).Value = ""

# Row 56  <=  case 'A 49032-2022'
$ws.Cells.Item(56, 1).Value = 'A 49032-2022'
$ws.Cells.Item(56, 2).Value = 44860.44083333333
$rngGQ = $ws.Range("G56:Q56")
$rngGQ.Value = @(3.7,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(56, 18).Value = ""

# Row 57  <=  case 'A 18121-2025'
$ws.Cells.Item(57, 1).Value = 'A 18121-2025'
$ws.Cells.Item(57, 2).Value = 45761
$rngGQ = $ws.Range("G57:Q57")
$rngGQ.Value = @(1.7,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(57, 18).Value = ""

# Row 58  <=  case 'A 43294-2025'
$ws.Cells.Item(58, 1).Value = 'A 43294-2025'
$ws.Cells.Item(58, 2).Value = 45910
$rngGQ = $ws.Range("G58:Q58")
$rngGQ.Value = @(2.6,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(58, 18).Value = ""

# Row 59  <=  case 'A 6335-2022'
$ws.Cells.Item(59, 1).Value = 'A 6335-2022'
$ws.Cells.Item(59, 2).Value = 44600
$rngGQ = $ws.Range("G59:Q59")
$rngGQ.Value = @(4,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(59, 18).Value = ""

# Row 60  <=  case 'A 20125-2024'
$ws.Cells.Item(60, 1).Value = 'A 20125-2024'
$ws.Cells.Item(60, 2).Value = 45434.55394675926
$rngGQ = $ws.Range("G60:Q60")
$rngGQ.Value = @(5.7,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(60, 18).Value = ""

# Row 61  <=  case 'A 10196-2025'
$ws.Cells.Item(61, 1).Value = 'A 10196-2025'
$ws.Cells.Item(61, 2).Value = 45719
$rngGQ = $ws.Range("G61:Q61")
$rngGQ.Value = @(0.5,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(61, 18).Value = ""

# Row 62  <=  case 'A 44112-2023'
$ws.Cells.Item(62, 1).Value = 'A 44112-2023'
$ws.Cells.Item(62, 2).Value = 45188.4418287037
$rngGQ = $ws.Range("G62:Q62")
$rngGQ.Value = @(0.5,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(62, 18).Value = ""

# Row 63  <=  case 'A 58970-2025'
$ws.Cells.Item(63, 1).Value = 'A 58970-2025'
$ws.Cells.Item(63, 2).Value = 45987.61269675926
$rngGQ = $ws.Range("G63:Q63")
$rngGQ.Value = @(1,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(63, 18).Value = ""

# Row 64  <=  case 'A 59230-2025'
$ws.Cells.Item(64, 1).Value = 'A 59230-2025'
$ws.Cells.Item(64, 2).Value = 45988.62502314815
$rngGQ = $ws.Range("G64:Q64")
$rngGQ.Value = @(0.6,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(64, 18).Value = ""

# Row 65  <=  case 'A 35911-2023'
$ws.Cells.Item(65, 1).Value = 'A 35911-2023'
$ws.Cells.Item(65, 2).Value = 45148.64078703704
$rngGQ = $ws.Range("G65:Q65")
$rngGQ.Value = @(0.9,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(65, 18).Value = ""

# Row 66  <=  case 'A 43466-2025'
$ws.Cells.Item(66, 1).Value = 'A 43466-2025'
$ws.Cells.Item(66, 2).Value = 45911
$rngGQ = $ws.Range("G66:Q66")
$rngGQ.Value = @(4.5,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(66, 18).Value = ""

# Row 67  <=  case 'A 11672-2024'
$ws.Cells.Item(67, 1).Value = 'A 11672-2024'
$ws.Cells.Item(67, 2).Value = 45373
$rngGQ = $ws.Range("G67:Q67")
$rngGQ.Value = @(13.9,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(67, 18).Value = ""

# Row 68  <=  case 'A 42588-2023'
$ws.Cells.Item(68, 1).Value = 'A 42588-2023'
$ws.Cells.Item(68, 2).Value = 45176
$rngGQ = $ws.Range("G68:Q68")
$rngGQ.Value = @(1.1,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(68, 18).Value = ""

# Row 69  <=  case 'A 1901-2025'
$ws.Cells.Item(69, 1).Value = 'A 1901-2025'
$ws.Cells.Item(69, 2).Value = 45671
$rngGQ = $ws.Range("G69:Q69")
$rngGQ.Value = @(4.2,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(69, 18).Value = ""

# Row 70  <=  case 'A 1920-2025'
$ws.Cells.Item(70, 1).Value = 'A 1920-2025'
$ws.Cells.Item(70, 2).Value = 45671
$rngGQ = $ws.Range("G70:Q70")
$rngGQ.Value = @(2.3,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(70, 18).Value = ""

# Row 71  <=  case 'A 43736-2023'
$ws.Cells.Item(71, 1).Value = 'A 43736-2023'
$ws.Cells.Item(71, 2).Value = 45182
$rngGQ = $ws.Range("G71:Q71")
$rngGQ.Value = @(1.2,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(71, 18).Value = ""

# Row 72  <=  case 'A 52664-2023'
$ws.Cells.Item(72, 1).Value = 'A 52664-2023'
$ws.Cells.Item(72, 2).Value = 45225
$rngGQ = $ws.Range("G72:Q72")
$rngGQ.Value = @(6.4,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(72, 18).Value = ""

# Row 73  <=  case 'A 3430-2026'
$ws.Cells.Item(73, 1).Value = 'A 3430-2026'
$ws.Cells.Item(73, 2).Value = 46042
$rngGQ = $ws.Range("G73:Q73")
$rngGQ.Value = @(5.7,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(73, 18).Value = ""

# Row 74  <=  case 'A 9716-2022'
$ws.Cells.Item(74, 1).Value = 'A 9716-2022'
$ws.Cells.Item(74, 2).Value = 44617
$rngGQ = $ws.Range("G74:Q74")
$rngGQ.Value = @(0.7,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(74, 18).Value = ""

# Row 75  <=  case 'A 43004-2025'
$ws.Cells.Item(75, 1).Value = 'A 43004-2025'
$ws.Cells.Item(75, 2).Value = 45909
$rngGQ = $ws.Range("G75:Q75")
$rngGQ.Value = @(1.6,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(75, 18).Value = ""

# Row 76  <=  case 'A 61045-2024'
$ws.Cells.Item(76, 1).Value = 'A 61045-2024'
$ws.Cells.Item(76, 2).Value = 45645.47678240741
$rngGQ = $ws.Range("G76:Q76")
$rngGQ.Value = @(5.7,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(76, 18).Value = ""

# Row 78  <=  case 'A 5068-2026'
$ws.Cells.Item(78, 1).Value = 'A 5068-2026'
$ws.Cells.Item(78, 2).Value = 46049.50232638889
$rngGQ = $ws.Range("G78:Q78")
$rngGQ.Value = @(3.5,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(78, 18).Value = ""

# Row 80  <=  case 'A 5190-2024'
$ws.Cells.Item(80, 1).Value = 'A 5190-2024'
$ws.Cells.Item(80, 2).Value = 45330
$rngGQ = $ws.Range("G80:Q80")
$rngGQ.Value = @(2.1,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(80, 18).Value = ""

# Row 81  <=  case 'A 17731-2025'
$ws.Cells.Item(81, 1).Value = 'A 17731-2025'
$ws.Cells.Item(81, 2).Value = 45758
$rngGQ = $ws.Range("G81:Q81")
$rngGQ.Value = @(0.8,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(81, 18).Value = ""

# Row 83  <=  case 'A 22792-2024'
$ws.Cells.Item(83, 1).Value = 'A 22792-2024'
$ws.Cells.Item(83, 2).Value = 45448.48983796296
$rngGQ = $ws.Range("G83:Q83")
$rngGQ.Value = @(0.5,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(83, 18).Value = ""

# Row 84  <=  case 'A 4964-2022'
$ws.Cells.Item(84, 1).Value = 'A 4964-2022'
$ws.Cells.Item(84, 2).Value = 44593
$rngGQ = $ws.Range("G84:Q84")
$rngGQ.Value = @(1.5,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(84, 18).Value = ""

# Row 85  <=  case 'A 11177-2023'
$ws.Cells.Item(85, 1).Value = 'A 11177-2023'
$ws.Cells.Item(85, 2).Value = 44987
$rngGQ = $ws.Range("G85:Q85")
$rngGQ.Value = @(1.2,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(85, 18).Value = ""

# Row 86  <=  case 'A 16215-2022'
$ws.Cells.Item(86, 1).Value = 'A 16215-2022'
$ws.Cells.Item(86, 2).Value = 44670
$rngGQ = $ws.Range("G86:Q86")
$rngGQ.Value = @(2.9,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(86, 18).Value = ""

# Row 88  <=  case 'A 48262-2021'
$ws.Cells.Item(88, 1).Value = 'A 48262-2021'
$ws.Cells.Item(88, 2).Value = 44449.66324074074
$rngGQ = $ws.Range("G88:Q88")
$rngGQ.Value = @(0.6,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(88, 18).Value = ""

# Row 89  <=  case 'A 8924-2023'
$ws.Cells.Item(89, 1).Value = 'A 8924-2023'
$ws.Cells.Item(89, 2).Value = 44974
$rngGQ = $ws.Range("G89:Q89")
$rngGQ.Value = @(6.8,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(89, 18).Value = ""

# Row 90  <=  case 'A 497-2023'
$ws.Cells.Item(90, 1).Value = 'A 497-2023'
$ws.Cells.Item(90, 2).Value = 44930.33540509259
$rngGQ = $ws.Range("G90:Q90")
$rngGQ.Value = @(0.8,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(90, 18).Value = ""

# Row 91  <=  case 'A 11633-2024'
$ws.Cells.Item(91, 1).Value = 'A 11633-2024'
$ws.Cells.Item(91, 2).Value = 45373.4740625
$rngGQ = $ws.Range("G91:Q91")
$rngGQ.Value = @(4.9,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(91, 18).Value = ""

# Row 92  <=  case 'A 11922-2023'
$ws.Cells.Item(92, 1).Value = 'A 11922-2023'
$ws.Cells.Item(92, 2).Value = 44993
$rngGQ = $ws.Range("G92:Q92")
$rngGQ.Value = @(0.9,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(92, 18).Value = ""

# Row 93  <=  case 'A 14546-2024'
$ws.Cells.Item(93, 1).Value = 'A 14546-2024'
$ws.Cells.Item(93, 2).Value = 45394
$rngGQ = $ws.Range("G93:Q93")
$rngGQ.Value = @(1,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(93, 18).Value = ""

# Row 94  <=  case 'A 32620-2023'
$ws.Cells.Item(94, 1).Value = 'A 32620-2023'
$ws.Cells.Item(94, 2).Value = 45111
$rngGQ = $ws.Range("G94:Q94")
$rngGQ.Value = @(1.4,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(94, 18).Value = ""

# Row 95  <=  case 'A 1592-2024'
$ws.Cells.Item(95, 1).Value = 'A 1592-2024'
$ws.Cells.Item(95, 2).Value = 45306.58594907408
$rngGQ = $ws.Range("G95:Q95")
$rngGQ.Value = @(2.2,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(95, 18).Value = ""

# Row 96  <=  case 'A 4149-2022'
$ws.Cells.Item(96, 1).Value = 'A 4149-2022'
$ws.Cells.Item(96, 2).Value = 44588
$rngGQ = $ws.Range("G96:Q96")
$rngGQ.Value = @(1.1,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(96, 18).Value = ""

# Row 97  <=  case 'A 15487-2025'
$ws.Cells.Item(97, 1).Value = 'A 15487-2025'
$ws.Cells.Item(97, 2).Value = 45747
$rngGQ = $ws.Range("G97:Q97")
$rngGQ.Value = @(0.9,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(97, 18).Value = ""

# Row 98  <=  case 'A 11101-2023'
$ws.Cells.Item(98, 1).Value = 'A 11101-2023'
$ws.Cells.Item(98, 2).Value = 44986
$rngGQ = $ws.Range("G98:Q98")
$rngGQ.Value = @(1.5,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(98, 18).Value = ""

# Row 100  <=  case 'A 13273-2024'
$ws.Cells.Item(100, 1).Value = 'A 13273-2024'
$ws.Cells.Item(100, 2).Value = 45386
$rngGQ = $ws.Range("G100:Q100")
$rngGQ.Value = @(12.4,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(100, 18).Value = ""

# Row 101  <=  case 'A 44196-2024'
$ws.Cells.Item(101, 1).Value = 'A 44196-2024'
$ws.Cells.Item(101, 2).Value = 45573.31803240741
$rngGQ = $ws.Range("G101:Q101")
$rngGQ.Value = @(1.8,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(101, 18).Value = ""

# Row 102  <=  case 'A 21507-2025'
$ws.Cells.Item(102, 1).Value = 'A 21507-2025'
$ws.Cells.Item(102, 2).Value = 45782.61987268519
$rngGQ = $ws.Range("G102:Q102")
$rngGQ.Value = @(1.7,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(102, 18).Value = ""

# Row 103  <=  case 'A 21556-2025'
$ws.Cells.Item(103, 1).Value = 'A 21556-2025'
$ws.Cells.Item(103, 2).Value = 45782.67826388889
$rngGQ = $ws.Range("G103:Q103")
$rngGQ.Value = @(1.3,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(103, 18).Value = ""

# Row 104  <=  case 'A 21541-2025'
$ws.Cells.Item(104, 1).Value = 'A 21541-2025'
$ws.Cells.Item(104, 2).Value = 45782.66538194445
$rngGQ = $ws.Range("G104:Q104")
$rngGQ.Value = @(3.1,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(104, 18).Value = ""

# Row 106  <=  case 'A 22203-2025'
$ws.Cells.Item(106, 1).Value = 'A 22203-2025'
$ws.Cells.Item(106, 2).Value = 45785.65381944444
$rngGQ = $ws.Range("G106:Q106")
$rngGQ.Value = @(2.8,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(106, 18).Value = ""

# Row 107  <=  case 'A 20047-2024'
$ws.Cells.Item(107, 1).Value = 'A 20047-2024'
$ws.Cells.Item(107, 2).Value = 45434.37376157408
$rngGQ = $ws.Range("G107:Q107")
$rngGQ.Value = @(2.3,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(107, 18).Value = ""

# Row 108  <=  case 'A 34808-2024'
$ws.Cells.Item(108, 1).Value = 'A 34808-2024'
$ws.Cells.Item(108, 2).Value = 45526.69388888889
$rngGQ = $ws.Range("G108:Q108")
$rngGQ.Value = @(3.6,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(108, 18).Value = ""

# Row 109  <=  case 'A 5167-2024'
$ws.Cells.Item(109, 1).Value = 'A 5167-2024'
$ws.Cells.Item(109, 2).Value = 45330
$rngGQ = $ws.Range("G109:Q109")
$rngGQ.Value = @(6.5,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(109, 18).Value = ""

# Row 110  <=  case 'A 23747-2025'
$ws.Cells.Item(110, 1).Value = 'A 23747-2025'
$ws.Cells.Item(110, 2).Value = 45793.44744212963
$rngGQ = $ws.Range("G110:Q110")
$rngGQ.Value = @(2.6,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(110, 18).Value = ""

# Row 111  <=  case 'A 24513-2025'
$ws.Cells.Item(111, 1).Value = 'A 24513-2025'
$ws.Cells.Item(111, 2).Value = 45798.40712962963
$rngGQ = $ws.Range("G111:Q111")
$rngGQ.Value = @(0.4,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(111, 18).Value = ""

# Row 112  <=  case 'A 27310-2023'
$ws.Cells.Item(112, 1).Value = 'A 27310-2023'
$ws.Cells.Item(112, 2).Value = 45096.65549768518
$rngGQ = $ws.Range("G112:Q112")
$rngGQ.Value = @(2.5,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(112, 18).Value = ""

# Row 113  <=  case 'A 26949-2025'
$ws.Cells.Item(113, 1).Value = 'A 26949-2025'
$ws.Cells.Item(113, 2).Value = 45811.43701388889
$rngGQ = $ws.Range("G113:Q113")
$rngGQ.Value = @(4.2,0,0,0,0,0,0,0,0,0,0)
$ws.Cells.Item(113, 18).Value = ""

# ---- 3) Column F ('Markagare') special-case: the 'Kyrkan' landowner note moves with its case ----
$ws.Cells.Item(55, 6).Value = ""
$ws.Cells.Item(48, 6).Value = 'Kyrkan'

# ---- 4) Refresh the HYPERLINK formulas (S,T,V,W,X,Y) for the case-lookup rows (2-13) ----
# These columns only exist for rows that carry a species/article note in column R.
# Row 2: case 'A 43326-2025'
$ws.Cells.Item(2, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/artfynd/A 43326-2025 artfynd.xlsx", "A 43326-2025")'
$ws.Cells.Item(2, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/kartor/A 43326-2025 karta.png", "A 43326-2025")'
$ws.Cells.Item(2, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/klagomål/A 43326-2025 FSC-klagomål.docx", "A 43326-2025")'
$ws.Cells.Item(2, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/klagomålsmail/A 43326-2025 FSC-klagomål mail.docx", "A 43326-2025")'
$ws.Cells.Item(2, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/tillsyn/A 43326-2025 tillsynsbegäran.docx", "A 43326-2025")'
$ws.Cells.Item(2, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/tillsynsmail/A 43326-2025 tillsynsbegäran mail.docx", "A 43326-2025")'

# Row 3: case 'A 39165-2024'
$ws.Cells.Item(3, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/artfynd/A 39165-2024 artfynd.xlsx", "A 39165-2024")'
$ws.Cells.Item(3, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/kartor/A 39165-2024 karta.png", "A 39165-2024")'
$ws.Cells.Item(3, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/klagomål/A 39165-2024 FSC-klagomål.docx", "A 39165-2024")'
$ws.Cells.Item(3, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/klagomålsmail/A 39165-2024 FSC-klagomål mail.docx", "A 39165-2024")'
$ws.Cells.Item(3, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/tillsyn/A 39165-2024 tillsynsbegäran.docx", "A 39165-2024")'
$ws.Cells.Item(3, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/tillsynsmail/A 39165-2024 tillsynsbegäran mail.docx", "A 39165-2024")'

# Row 4: case 'A 38522-2025'
$ws.Cells.Item(4, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/artfynd/A 38522-2025 artfynd.xlsx", "A 38522-2025")'
$ws.Cells.Item(4, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/kartor/A 38522-2025 karta.png", "A 38522-2025")'
$ws.Cells.Item(4, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/klagomål/A 38522-2025 FSC-klagomål.docx", "A 38522-2025")'
$ws.Cells.Item(4, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/klagomålsmail/A 38522-2025 FSC-klagomål mail.docx", "A 38522-2025")'
$ws.Cells.Item(4, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/tillsyn/A 38522-2025 tillsynsbegäran.docx", "A 38522-2025")'
$ws.Cells.Item(4, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/tillsynsmail/A 38522-2025 tillsynsbegäran mail.docx", "A 38522-2025")'

# Row 5: case 'A 61064-2024'
$ws.Cells.Item(5, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/artfynd/A 61064-2024 artfynd.xlsx", "A 61064-2024")'
$ws.Cells.Item(5, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/kartor/A 61064-2024 karta.png", "A 61064-2024")'
$ws.Cells.Item(5, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/klagomål/A 61064-2024 FSC-klagomål.docx", "A 61064-2024")'
$ws.Cells.Item(5, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/klagomålsmail/A 61064-2024 FSC-klagomål mail.docx", "A 61064-2024")'
$ws.Cells.Item(5, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/tillsyn/A 61064-2024 tillsynsbegäran.docx", "A 61064-2024")'
$ws.Cells.Item(5, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/tillsynsmail/A 61064-2024 tillsynsbegäran mail.docx", "A 61064-2024")'

# Row 6: case 'A 26855-2022'
$ws.Cells.Item(6, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/artfynd/A 26855-2022 artfynd.xlsx", "A 26855-2022")'
$ws.Cells.Item(6, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/kartor/A 26855-2022 karta.png", "A 26855-2022")'
$ws.Cells.Item(6, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/klagomål/A 26855-2022 FSC-klagomål.docx", "A 26855-2022")'
$ws.Cells.Item(6, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/klagomålsmail/A 26855-2022 FSC-klagomål mail.docx", "A 26855-2022")'
$ws.Cells.Item(6, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/tillsyn/A 26855-2022 tillsynsbegäran.docx", "A 26855-2022")'
$ws.Cells.Item(6, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/tillsynsmail/A 26855-2022 tillsynsbegäran mail.docx", "A 26855-2022")'

# Row 7: case 'A 15456-2024'
$ws.Cells.Item(7, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/artfynd/A 15456-2024 artfynd.xlsx", "A 15456-2024")'
$ws.Cells.Item(7, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/kartor/A 15456-2024 karta.png", "A 15456-2024")'
$ws.Cells.Item(7, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/klagomål/A 15456-2024 FSC-klagomål.docx", "A 15456-2024")'
$ws.Cells.Item(7, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/klagomålsmail/A 15456-2024 FSC-klagomål mail.docx", "A 15456-2024")'
$ws.Cells.Item(7, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/tillsyn/A 15456-2024 tillsynsbegäran.docx", "A 15456-2024")'
$ws.Cells.Item(7, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/tillsynsmail/A 15456-2024 tillsynsbegäran mail.docx", "A 15456-2024")'

# Row 8: case 'A 15475-2024'
$ws.Cells.Item(8, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/artfynd/A 15475-2024 artfynd.xlsx", "A 15475-2024")'
$ws.Cells.Item(8, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/kartor/A 15475-2024 karta.png", "A 15475-2024")'
$ws.Cells.Item(8, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/klagomål/A 15475-2024 FSC-klagomål.docx", "A 15475-2024")'
$ws.Cells.Item(8, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/klagomålsmail/A 15475-2024 FSC-klagomål mail.docx", "A 15475-2024")'
$ws.Cells.Item(8, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/tillsyn/A 15475-2024 tillsynsbegäran.docx", "A 15475-2024")'
$ws.Cells.Item(8, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/tillsynsmail/A 15475-2024 tillsynsbegäran mail.docx", "A 15475-2024")'

# Row 9: case 'A 54424-2023'
$ws.Cells.Item(9, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/artfynd/A 54424-2023 artfynd.xlsx", "A 54424-2023")'
$ws.Cells.Item(9, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/kartor/A 54424-2023 karta.png", "A 54424-2023")'
$ws.Cells.Item(9, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/klagomål/A 54424-2023 FSC-klagomål.docx", "A 54424-2023")'
$ws.Cells.Item(9, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/klagomålsmail/A 54424-2023 FSC-klagomål mail.docx", "A 54424-2023")'
$ws.Cells.Item(9, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/tillsyn/A 54424-2023 tillsynsbegäran.docx", "A 54424-2023")'
$ws.Cells.Item(9, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/tillsynsmail/A 54424-2023 tillsynsbegäran mail.docx", "A 54424-2023")'

# Row 10: case 'A 1577-2024'
$ws.Cells.Item(10, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/artfynd/A 1577-2024 artfynd.xlsx", "A 1577-2024")'
$ws.Cells.Item(10, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/kartor/A 1577-2024 karta.png", "A 1577-2024")'
$ws.Cells.Item(10, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/klagomål/A 1577-2024 FSC-klagomål.docx", "A 1577-2024")'
$ws.Cells.Item(10, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/klagomålsmail/A 1577-2024 FSC-klagomål mail.docx", "A 1577-2024")'
$ws.Cells.Item(10, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/tillsyn/A 1577-2024 tillsynsbegäran.docx", "A 1577-2024")'
$ws.Cells.Item(10, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/tillsynsmail/A 1577-2024 tillsynsbegäran mail.docx", "A 1577-2024")'

# Row 11: case 'A 54127-2025'
$ws.Cells.Item(11, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/artfynd/A 54127-2025 artfynd.xlsx", "A 54127-2025")'
$ws.Cells.Item(11, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/kartor/A 54127-2025 karta.png", "A 54127-2025")'
$ws.Cells.Item(11, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/klagomål/A 54127-2025 FSC-klagomål.docx", "A 54127-2025")'
$ws.Cells.Item(11, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/klagomålsmail/A 54127-2025 FSC-klagomål mail.docx", "A 54127-2025")'
$ws.Cells.Item(11, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/tillsyn/A 54127-2025 tillsynsbegäran.docx", "A 54127-2025")'
$ws.Cells.Item(11, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/tillsynsmail/A 54127-2025 tillsynsbegäran mail.docx", "A 54127-2025")'

# Row 12: case 'A 11170-2023'
$ws.Cells.Item(12, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/artfynd/A 11170-2023 artfynd.xlsx", "A 11170-2023")'
$ws.Cells.Item(12, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/kartor/A 11170-2023 karta.png", "A 11170-2023")'
$ws.Cells.Item(12, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/klagomål/A 11170-2023 FSC-klagomål.docx", "A 11170-2023")'
$ws.Cells.Item(12, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/klagomålsmail/A 11170-2023 FSC-klagomål mail.docx", "A 11170-2023")'
$ws.Cells.Item(12, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/tillsyn/A 11170-2023 tillsynsbegäran.docx", "A 11170-2023")'
$ws.Cells.Item(12, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/tillsynsmail/A 11170-2023 tillsynsbegäran mail.docx", "A 11170-2023")'

# Row 13: case 'A 8169-2024'
$ws.Cells.Item(13, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/artfynd/A 8169-2024 artfynd.xlsx", "A 8169-2024")'
$ws.Cells.Item(13, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/kartor/A 8169-2024 karta.png", "A 8169-2024")'
$ws.Cells.Item(13, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/klagomål/A 8169-2024 FSC-klagomål.docx", "A 8169-2024")'
$ws.Cells.Item(13, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/klagomålsmail/A 8169-2024 FSC-klagomål mail.docx", "A 8169-2024")'
$ws.Cells.Item(13, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/tillsyn/A 8169-2024 tillsynsbegäran.docx", "A 8169-2024")'
$ws.Cells.Item(13, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_1214/tillsynsmail/A 8169-2024 tillsynsbegäran mail.docx", "A 8169-2024")'
